$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 20740
$ws.Range("J87").Value = 20740
$ws.Range("L87").Value = 20740
$ws.Range("N87").Value = -23236

$ws.Range("H90").Value = 20740
$ws.Range("J90").Value = 20740
$ws.Range("L90").Value = 62220
$ws.Range("N90").Value = -74700

$ws.Range("H98").Value = 1853.1
$ws.Range("I98").Value = 1198.6111
$ws.Range("K98").Value = 1198.6111
$ws.Range("M98").Value = 299.3888999999999

$ws.Range("H122").Value = 1853.1
$ws.Range("I122").Value = 1198.6111
$ws.Range("K122").Value = 3595.8333
$ws.Range("M122").Value = -1145.8333

$ws.Range("H129").Value = 1086.5
$ws.Range("I129").Value = 350
$ws.Range("J129").Value = 1191.7142
$ws.Range("K129").Value = 1050
$ws.Range("L129").Value = 3575.1426
$ws.Range("M129").Value = 3950
$ws.Range("N129").Value = -13575.1426

$ws.Range("H137").Value = 3063.7727
$ws.Range("I137").Value = 3558.24
$ws.Range("K137").Value = 10674.72
$ws.Range("M137").Value = -8124.719999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1380.0588
$ws.Range("I2").Value = 1358.5385
$ws.Range("J2").Value = 1450
$ws.Range("K2").Value = 1358.5385
$ws.Range("L2").Value = 1450
$ws.Range("M2").Value = -1245.5385
$ws.Range("N2").Value = -1676

$ws.Range("H37").Value = 30579.637
$ws.Range("J37").Value = 34038
$ws.Range("L37").Value = 34038
$ws.Range("N37").Value = -34584

$ws.Range("H45").Value = 1549.4584
$ws.Range("I45").Value = 1544.2273
$ws.Range("J45").Value = 1607
$ws.Range("K45").Value = 1544.2273
$ws.Range("L45").Value = 1607
$ws.Range("M45").Value = -1167.2273
$ws.Range("N45").Value = -2361

$ws.Range("H116").Value = 1380.0588
$ws.Range("I116").Value = 1358.5385
$ws.Range("J116").Value = 1450
$ws.Range("K116").Value = 1358.5385
$ws.Range("L116").Value = 1450
$ws.Range("M116").Value = 935.4614999999999
$ws.Range("N116").Value = -6038

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1380.0588
$ws.Range("I3").Value = 1358.5385
$ws.Range("J3").Value = 1450
$ws.Range("K3").Value = 1358.5385
$ws.Range("L3").Value = 1450
$ws.Range("M3").Value = -1244.5385
$ws.Range("N3").Value = -1678

$ws.Range("H64").Value = 1000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 1000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 1000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -1450

$ws.Range("H67").Value = 1000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 1000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 1000
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -2560

$ws.Range("H107").Value = 626.3333
$ws.Range("I107").Value = 633.2727
$ws.Range("J107").Value = 550
$ws.Range("K107").Value = 633.2727
$ws.Range("L107").Value = 550
$ws.Range("M107").Value = 1286.7273
$ws.Range("N107").Value = -4390

$ws.Range("H134").Value = 1629.7333
$ws.Range("I134").Value = 1290.75
$ws.Range("K134").Value = 3872.25
$ws.Range("M134").Value = -1337.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1057.6923
$ws.Range("I22").Value = 1373.3334
$ws.Range("J22").Value = 347.5
$ws.Range("K22").Value = 1373.3334
$ws.Range("L22").Value = 347.5
$ws.Range("M22").Value = -1023.3334
$ws.Range("N22").Value = -1047.5

$ws.Range("H99").Value = 1726.0625
$ws.Range("I99").Value = 1662.2222
$ws.Range("J99").Value = 1808.1428
$ws.Range("K99").Value = 1662.2222
$ws.Range("L99").Value = 1808.1428
$ws.Range("M99").Value = -164.2221999999999
$ws.Range("N99").Value = -4804.1428

$ws.Range("H126").Value = 1726.0625
$ws.Range("I126").Value = 1662.2222
$ws.Range("J126").Value = 1808.1428
$ws.Range("K126").Value = 4986.6666
$ws.Range("L126").Value = 5424.428400000001
$ws.Range("M126").Value = -2516.6666
$ws.Range("N126").Value = -10364.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 1065.4445
$ws.Range("I44").Value = 422.25
$ws.Range("J44").Value = 1580
$ws.Range("K44").Value = 1266.75
$ws.Range("L44").Value = 4740
$ws.Range("M44").Value = -868.75
$ws.Range("N44").Value = -5536

$ws.Range("H98").Value = 770523.25
$ws.Range("I98").Value = 618.4545000000001
$ws.Range("K98").Value = 1855.3635
$ws.Range("M98").Value = -357.3635000000002

$ws.Range("H118").Value = 2346.8
$ws.Range("I118").Value = 762.8
$ws.Range("J118").Value = 2610.8
$ws.Range("K118").Value = 2288.4
$ws.Range("L118").Value = 7832.400000000001
$ws.Range("M118").Value = -1045.4
$ws.Range("N118").Value = -10318.4

$ws.Range("H131").Value = 1327027.4
$ws.Range("I131").Value = 360.41666
$ws.Range("J131").Value = 1673114.5
$ws.Range("K131").Value = 1081.24998
$ws.Range("L131").Value = 5019343.5
$ws.Range("M131").Value = 3958.75002
$ws.Range("N131").Value = -5029423.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 80000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 80000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 80000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -80302

$ws.Range("H70").Value = 4624.7188
$ws.Range("I70").Value = 4325.68
$ws.Range("J70").Value = 5692.7144
$ws.Range("K70").Value = 4325.68
$ws.Range("L70").Value = 5692.7144
$ws.Range("M70").Value = -4055.68
$ws.Range("N70").Value = -6232.7144

$ws.Range("H73").Value = 4624.7188
$ws.Range("I73").Value = 4325.68
$ws.Range("J73").Value = 5692.7144
$ws.Range("K73").Value = 4325.68
$ws.Range("L73").Value = 5692.7144
$ws.Range("M73").Value = -3389.68
$ws.Range("N73").Value = -7564.7144

$ws.Range("H102").Value = 1316.091
$ws.Range("I102").Value = 1187.619
$ws.Range("J102").Value = 4014
$ws.Range("K102").Value = 1187.619
$ws.Range("L102").Value = 4014
$ws.Range("M102").Value = 434.3810000000001
$ws.Range("N102").Value = -7258

$ws.Range("H126").Value = 101133.2
$ws.Range("I126").Value = 250848
$ws.Range("J126").Value = 1323.3334
$ws.Range("K126").Value = 752544
$ws.Range("L126").Value = 3970.0002
$ws.Range("M126").Value = -750074
$ws.Range("N126").Value = -8910.0002

$ws.Range("H132").Value = 2995.8667
$ws.Range("I132").Value = 2654
$ws.Range("J132").Value = 3679.6
$ws.Range("K132").Value = 7962
$ws.Range("L132").Value = 11038.8
$ws.Range("M132").Value = -5432
$ws.Range("N132").Value = -16098.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2100.7144
$ws.Range("I7").Value = 1480
$ws.Range("J7").Value = 3652.5
$ws.Range("K7").Value = 1480
$ws.Range("L7").Value = 3652.5
$ws.Range("M7").Value = -1368
$ws.Range("N7").Value = -3876.5

$ws.Range("H22").Value = 921.53845
$ws.Range("I22").Value = 876
$ws.Range("J22").Value = 950
$ws.Range("K22").Value = 876
$ws.Range("L22").Value = 950
$ws.Range("M22").Value = -581
$ws.Range("N22").Value = -1540

$ws.Range("H27").Value = 921.53845
$ws.Range("I27").Value = 876
$ws.Range("J27").Value = 950
$ws.Range("K27").Value = 876
$ws.Range("L27").Value = 950
$ws.Range("M27").Value = -769
$ws.Range("N27").Value = -1164

$ws.Range("H40").Value = 3655.182
$ws.Range("I40").Value = 3220.7
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 3220.7
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -3084.7
$ws.Range("N40").Value = -8272

$ws.Range("H46").Value = 626.6667
$ws.Range("I46").Value = 626.6667
$ws.Range("K46").Value = 626.6667
$ws.Range("M46").Value = -438.6667

$ws.Range("H55").Value = 546.5
$ws.Range("J55").Value = 750
$ws.Range("L55").Value = 750
$ws.Range("N55").Value = -1096

$ws.Range("H58").Value = 1100
$ws.Range("I58").Value = 1100
$ws.Range("K58").Value = 1100
$ws.Range("M58").Value = -840

$ws.Range("H126").Value = 2100.7144
$ws.Range("I126").Value = 1480
$ws.Range("J126").Value = 3652.5
$ws.Range("K126").Value = 4440
$ws.Range("L126").Value = 10957.5
$ws.Range("M126").Value = -1970
$ws.Range("N126").Value = -15897.5

$ws.Range("H132").Value = 6262.074
$ws.Range("I132").Value = 6792.7896
$ws.Range("J132").Value = 5001.625
$ws.Range("K132").Value = 20378.3688
$ws.Range("L132").Value = 15004.875
$ws.Range("M132").Value = -17848.3688
$ws.Range("N132").Value = -20064.875

$ws.Range("H133").Value = 28035.2
$ws.Range("J133").Value = 28035.2
$ws.Range("L133").Value = 28035.2
$ws.Range("N133").Value = -33095.2

$ws.Range("H136").Value = 2660.776
$ws.Range("I136").Value = 2118.9062
$ws.Range("J136").Value = 3327.6924
$ws.Range("K136").Value = 6356.7186
$ws.Range("L136").Value = 9983.0772
$ws.Range("M136").Value = -3806.7186
$ws.Range("N136").Value = -15083.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 556342.75
$ws.Range("I126").Value = 625741.2
$ws.Range("J126").Value = 1155
$ws.Range("K126").Value = 1877223.6
$ws.Range("L126").Value = 3465
$ws.Range("M126").Value = -1874753.6
$ws.Range("N126").Value = -8405

$ws.Range("H132").Value = 1138.3615
$ws.Range("I132").Value = 936.3582
$ws.Range("J132").Value = 1984.25
$ws.Range("K132").Value = 2809.0746
$ws.Range("L132").Value = 5952.75
$ws.Range("M132").Value = -279.0745999999999
$ws.Range("N132").Value = -11012.75
